$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.727.30'
$ws.Range("E2").Value = '  +2.07%  '

$ws.Range("D3").Value = '2.342.91'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.82%  '

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.574'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.09%  '

$ws.Range("D9").Value = '2.340.49'
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("E10").Value = '  +1.22%  '

$ws.Range("E11").Value = '  +0.98%  '

$ws.Range("E12").Value = '  +3.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.359'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.70%  '

$ws.Range("D14").Value = '2.759.53'
$ws.Range("E14").Value = '  +0.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.41%  '

$ws.Range("D16").Value = '57.716.79'
$ws.Range("E16").Value = '  +2.11%  '

$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").Value = '2.352.06'
$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '333.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("E21").Value = '  +1.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.45%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.21%  '

$ws.Range("E25").Value = '  +2.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("D31").Value = '0.0₃0729'
$ws.Range("E31").Value = '  +1.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("E33").Value = '  +16.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.30%  '

$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.61'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.48%  '

$ws.Range("E42").Value = '  +1.21%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.78%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '284.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.25%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0947'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.00%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.57%  '

$ws.Range("E47").Value = '  +1.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.561'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("B49").Value = 'Polygon'
$ws.Range("C49").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.385'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.22%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0217'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.06%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.48%  '
